# fri26.xlsx — "Added most of Paden's code"
#
# The sheet originally stored only the lower-triangular half (incl. the
# zero diagonal) of the fri26 TSP symmetric edge-weight matrix in A8:Z33
# (one triangular row per worksheet row, growing from 1 cell in row 8 to
# 26 cells in row 33). The commit fills in the rest of each row so every
# row A8:Z33 holds the full, mirrored 26x26 distance matrix (value[i][j]
# == value[j][i]), and leaves the active selection on the (now larger)
# working area below the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full symmetric 26x26 fri26 edge-weight matrix (mirror of the existing
# lower triangle already present in A8:Z33).
$m = New-Object 'object[,]' 26,26

$row0 = 0,83,93,129,133,139,151,169,135,114,110,98,99,95,81,152,159,181,172,185,147,157,185,220,127,181
for ($j = 0; $j -lt 26; $j++) { $m[0,$j] = $row0[$j] }
$row1 = 83,0,40,53,62,64,91,116,93,84,95,98,89,68,67,127,156,175,152,165,160,180,223,268,179,197
for ($j = 0; $j -lt 26; $j++) { $m[1,$j] = $row1[$j] }
$row2 = 93,40,0,42,42,49,59,81,54,44,58,64,54,31,36,86,117,135,112,125,124,147,193,241,157,161
for ($j = 0; $j -lt 26; $j++) { $m[2,$j] = $row2[$j] }
$row3 = 129,53,42,0,11,11,46,72,65,70,88,100,89,66,76,102,142,156,127,139,155,180,228,278,197,190
for ($j = 0; $j -lt 26; $j++) { $m[3,$j] = $row3[$j] }
$row4 = 133,62,42,11,0,9,35,61,55,62,82,95,84,62,74,93,133,146,117,128,148,173,222,272,194,182
for ($j = 0; $j -lt 26; $j++) { $m[4,$j] = $row4[$j] }
$row5 = 139,64,49,11,9,0,39,65,63,71,90,103,92,71,82,100,141,153,124,135,156,181,230,280,202,190
for ($j = 0; $j -lt 26; $j++) { $m[5,$j] = $row5[$j] }
$row6 = 151,91,59,46,35,39,0,26,34,52,71,88,77,63,78,66,110,119,88,98,130,156,206,257,188,160
for ($j = 0; $j -lt 26; $j++) { $m[6,$j] = $row6[$j] }
$row7 = 169,116,81,72,61,65,26,0,37,59,75,92,83,76,91,54,98,103,70,78,122,148,198,250,188,148
for ($j = 0; $j -lt 26; $j++) { $m[7,$j] = $row7[$j] }
$row8 = 135,93,54,65,55,63,34,37,0,22,39,56,47,40,55,37,78,91,62,74,96,122,172,223,155,128
for ($j = 0; $j -lt 26; $j++) { $m[8,$j] = $row8[$j] }
$row9 = 114,84,44,70,62,71,52,59,22,0,20,36,26,20,34,43,74,91,68,82,86,111,160,210,136,121
for ($j = 0; $j -lt 26; $j++) { $m[9,$j] = $row9[$j] }
$row10 = 110,95,58,88,82,90,71,75,39,20,0,18,11,27,32,42,61,80,64,77,68,92,140,190,116,103
for ($j = 0; $j -lt 26; $j++) { $m[10,$j] = $row10[$j] }
$row11 = 98,98,64,100,95,103,88,92,56,36,18,0,11,34,31,56,63,85,75,87,62,83,129,178,100,99
for ($j = 0; $j -lt 26; $j++) { $m[11,$j] = $row11[$j] }
$row12 = 99,89,54,89,84,92,77,83,47,26,11,11,0,23,24,53,68,89,74,87,71,93,140,189,111,107
for ($j = 0; $j -lt 26; $j++) { $m[12,$j] = $row12[$j] }
$row13 = 95,68,31,66,62,71,63,76,40,20,27,34,23,0,15,62,87,106,87,100,93,116,163,212,132,130
for ($j = 0; $j -lt 26; $j++) { $m[13,$j] = $row13[$j] }
$row14 = 81,67,36,76,74,82,78,91,55,34,32,31,24,15,0,73,92,112,96,109,93,113,158,205,122,130
for ($j = 0; $j -lt 26; $j++) { $m[14,$j] = $row14[$j] }
$row15 = 152,127,86,102,93,100,66,54,37,43,42,56,53,62,73,0,44,54,26,39,68,94,144,196,139,95
for ($j = 0; $j -lt 26; $j++) { $m[15,$j] = $row15[$j] }
$row16 = 159,156,117,142,133,141,110,98,78,74,61,63,68,87,92,44,0,22,34,38,30,53,102,154,109,51
for ($j = 0; $j -lt 26; $j++) { $m[16,$j] = $row16[$j] }
$row17 = 181,175,135,156,146,153,119,103,91,91,80,85,89,106,112,54,22,0,33,29,46,64,107,157,125,51
for ($j = 0; $j -lt 26; $j++) { $m[17,$j] = $row17[$j] }
$row18 = 172,152,112,127,117,124,88,70,62,68,64,75,74,87,96,26,34,33,0,13,63,87,135,186,141,81
for ($j = 0; $j -lt 26; $j++) { $m[18,$j] = $row18[$j] }
$row19 = 185,165,125,139,128,135,98,78,74,82,77,87,87,100,109,39,38,29,13,0,68,90,136,186,148,79
for ($j = 0; $j -lt 26; $j++) { $m[19,$j] = $row19[$j] }
$row20 = 147,160,124,155,148,156,130,122,96,86,68,62,71,93,93,68,30,46,63,68,0,26,77,128,80,37
for ($j = 0; $j -lt 26; $j++) { $m[20,$j] = $row20[$j] }
$row21 = 157,180,147,180,173,181,156,148,122,111,92,83,93,116,113,94,53,64,87,90,26,0,50,102,65,27
for ($j = 0; $j -lt 26; $j++) { $m[21,$j] = $row21[$j] }
$row22 = 185,223,193,228,222,230,206,198,172,160,140,129,140,163,158,144,102,107,135,136,77,50,0,51,64,58
for ($j = 0; $j -lt 26; $j++) { $m[22,$j] = $row22[$j] }
$row23 = 220,268,241,278,272,280,257,250,223,210,190,178,189,212,205,196,154,157,186,186,128,102,51,0,93,107
for ($j = 0; $j -lt 26; $j++) { $m[23,$j] = $row23[$j] }
$row24 = 127,179,157,197,194,202,188,188,155,136,116,100,111,132,122,139,109,125,141,148,80,65,64,93,0,90
for ($j = 0; $j -lt 26; $j++) { $m[24,$j] = $row24[$j] }
$row25 = 181,197,161,190,182,190,160,148,128,121,103,99,107,130,130,95,51,51,81,79,37,27,58,107,90,0
for ($j = 0; $j -lt 26; $j++) { $m[25,$j] = $row25[$j] }

# Write the whole matrix back in one shot (rows 8-33, columns A-Z).
$ws.Range("A8:Z33").Value = $m

# Matches the author selecting the scratch area below the table before saving.
$ws.Range("B35:AA60").Select()

